$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.421.51"
$ws.Range("E2").Value = "  +4.22%  "

$ws.Range("D3").Value = "3.605.50"
$ws.Range("E3").Value = "  +3.88%  "

$ws.Range("E4").Value = "  +0.30%  "

$c = $ws.Range("D5")
$c.Value = "'628.50"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +3.95%  "

$c = $ws.Range("D6")
$c.Value = "'158.51"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +6.38%  "

$ws.Range("D7").Value = "3.605.78"
$ws.Range("E7").Value = "  +3.92%  "

$ws.Range("E8").Value = "  +0.08%  "

$c = $ws.Range("D9")
$c.Value = "'0.496"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +3.37%  "

$ws.Range("E10").Value = "  +8.00%  "

$c = $ws.Range("D11")
$c.Value = "'7.39"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +7.27%  "

$c = $ws.Range("D12")
$c.Value = "'0.440"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +4.43%  "

$ws.Range("E13").Value = "  +4.97%  "

$c = $ws.Range("D14")
$c.Value = "'33.53"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +6.73%  "

$ws.Range("D15").Value = "4.226.23"
$ws.Range("E15").Value = "  +4.29%  "

$ws.Range("D16").Value = "69.673.33"
$ws.Range("E16").Value = "  +4.76%  "

$ws.Range("D17").Value = "3.601.73"
$ws.Range("E17").Value = "  +4.21%  "

$ws.Range("E18").Value = "  +0.53%  "

$c = $ws.Range("D19")
$c.Value = "'6.70"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +5.41%  "

$c = $ws.Range("D20")
$c.Value = "'16.14"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +7.61%  "

$c = $ws.Range("D21")
$c.Value = "'10.24"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +13.42%  "

$c = $ws.Range("D22")
$c.Value = "'462.50"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +4.42%  "

$c = $ws.Range("D23")
$c.Value = "'0.646"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +3.92%  "

$c = $ws.Range("D24")
$c.Value = "'78.69"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +2.07%  "

$c = $ws.Range("D25")
$c.Value = "'0.0000135"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +9.62%  "

$c = $ws.Range("D26")
$c.Value = "'10.68"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +6.23%  "

$ws.Range("D27").Value = "3.758.71"
$ws.Range("E27").Value = "  +4.26%  "

$ws.Range("E28").Value = "  +0.09%  "

$c = $ws.Range("D29")
$c.Value = "'9.28"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +12.85%  "

$ws.Range("E30").Value = "  +5.04%  "

$c = $ws.Range("D31")
$c.Value = "'1.72"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +10.08%  "

$c = $ws.Range("D32")
$c.Value = "'0.176"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +11.42%  "

$c = $ws.Range("D33")
$c.Value = "'6.55"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +6.89%  "

$ws.Range("E34").Value = "  +0.13%  "

$c = $ws.Range("D35")
$c.Value = "'26.51"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +3.84%  "

$ws.Range("E36").Value = "  +4.88%  "

$ws.Range("D37").Value = "3.611.69"
$ws.Range("E37").Value = "  +4.53%  "

$c = $ws.Range("D38")
$c.Value = "'8.49"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +7.22%  "

$ws.Range("E39").Value = "  +11.24%  "

$ws.Range("E40").Value = "  +0.00%  "

$ws.Range("B41").Value = "Monero"
$ws.Range("C41").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$c = $ws.Range("D41")
$c.Value = "'180.01"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +3.91%  "

$ws.Range("B42").Value = "Hedera"
$ws.Range("C42").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$c = $ws.Range("D42")
$c.Value = "'0.0925"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +7.20%  "

$c = $ws.Range("D43")
$c.Value = "'1.00"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +0.22%  "

$c = $ws.Range("D44")
$c.Value = "'5.66"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +3.40%  "

$c = $ws.Range("D45")
$c.Value = "'31.95"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +21.33%  "

$c = $ws.Range("D46")
$c.Value = "'0.913"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +3.83%  "

$c = $ws.Range("D47")
$c.Value = "'1.39"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +12.79%  "

$ws.Range("B48").Value = "dogwifhat"
$ws.Range("C48").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$c = $ws.Range("D48")
$c.Value = "'2.75"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +9.70%  "

$ws.Range("B49").Value = "OKB"
$ws.Range("C49").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$c = $ws.Range("D49")
$c.Value = "'45.92"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +1.45%  "

$c = $ws.Range("D50")
$c.Value = "'7.81"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +3.59%  "

$c = $ws.Range("D51")
$c.Value = "'0.268"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +8.74%  "

Write-Host "done"